$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.413.82"
$ws.Range("E2").Value = "  +4.35%  "
$ws.Range("D3").Value = "2.350.93"
$ws.Range("E3").Value = "  +2.91%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "547.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.47"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.50%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.589"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.39%  "
$ws.Range("D9").Value = "2.346.44"
$ws.Range("E9").Value = "  +2.82%  "
$ws.Range("E10").Value = "  +1.91%  "
$ws.Range("E11").Value = "  +1.78%  "
$ws.Range("E12").Value = "  +1.24%  "
$ws.Range("E13").Value = "  +1.72%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.99"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.34%  "
$ws.Range("D15").Value = "2.770.63"
$ws.Range("E15").Value = "  +2.87%  "
$ws.Range("D16").Value = "60.421.73"
$ws.Range("E16").Value = "  +4.46%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000133"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.14%  "
$ws.Range("D18").Value = "2.355.86"
$ws.Range("E18").Value = "  +1.26%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.70"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.27%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.83%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.84"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +7.70%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "314.86"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.12%  "
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.21"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.51%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.171"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.64%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.96"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.27%  "
$ws.Range("E28").Value = "  +5.18%  "
$ws.Range("E29").Value = "  +3.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "171.81"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.81%  "
$ws.Range("B31").Value = "SuiNetwork"
$ws.Range("C31").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.15"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +10.40%  "
$ws.Range("B32").Value = "PEPE"
$ws.Range("C32").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D32").Value = "0.0₃0731"
$ws.Range("E32").Value = "  +2.24%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.92"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.41"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +14.77%  "
$ws.Range("E35").Value = "  +1.25%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.03"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.03%  "
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("E38").Value = "  +0.16%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.16"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.38%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "315.60"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +10.11%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "38.20"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.27%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.53"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.55%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "142.24"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.62%  "
$ws.Range("E44").Value = "  +1.91%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0954"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.10%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.31"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.98%  "
$ws.Range("E47").Value = "  +1.08%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.561"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.74%  "
$ws.Range("E49").Value = "  +2.38%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "11.05"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.95%  "
$ws.Range("D51").Value = "0.0₆0208"
$ws.Range("E51").Value = "  +7.79%  "
